# Labels.xlsx edit — per commit: "aip enhancements" (rebrand Mobimo -> Alya
# Consulting, fix Oeffentlich -> Öffentlich umlaut, and swap the EN/DE name
# & comment columns so English comes first).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture current A-D values for every populated row (1..12) ---
$rows = 1..12
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @(
        $ws.Cells.Item($r,1).Value2,
        $ws.Cells.Item($r,2).Value2,
        $ws.Cells.Item($r,3).Value2,
        $ws.Cells.Item($r,4).Value2
    )
}

function Rebrand([string]$s) {
    if ($null -eq $s) { return $s }
    $s = $s.Replace("Mobimo", "Alya Consulting")
    $s = $s.Replace("Oeffentlich", "Öffentlich")
    return $s
}

# --- 2. Write back with columns A<->C and B<->D swapped, applying the
#        text fixups to every cell touched ---
foreach ($r in $rows) {
    $a = $orig[$r][0]
    $b = $orig[$r][1]
    $c = $orig[$r][2]
    $d = $orig[$r][3]

    $newA = Rebrand($c)
    $newB = Rebrand($d)
    $newC = Rebrand($a)
    $newD = Rebrand($b)

    if ($null -ne $a -or $null -ne $c) { $ws.Cells.Item($r,1).Value = $newA }
    if ($null -ne $b -or $null -ne $d) { $ws.Cells.Item($r,2).Value = $newB }
    if ($null -ne $c -or $null -ne $a) { $ws.Cells.Item($r,3).Value = $newC }
    if ($null -ne $d -or $null -ne $b) { $ws.Cells.Item($r,4).Value = $newD }
}

# --- 3. Column A/C used to hold the (EN width=14.89 / DE width=14.62)
#        bestFit widths; now that the content has swapped, swap the
#        widths too ---
$colAWidth = $ws.Columns.Item(1).ColumnWidth
$colCWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(1).ColumnWidth = $colCWidth
$ws.Columns.Item(3).ColumnWidth = $colAWidth
